$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.124.29'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '2.759.63'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  +0.05%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.59'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -2.35%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.08'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -3.17%  '
$ws.Range('E9').Value = '  -4.13%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.88'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -13.59%  '
$ws.Range('E11').Value = '  +3.34%  '
$ws.Range('E12').Value = '  -2.88%  '
$ws.Range('D13').Value = '3.248.87'
$ws.Range('E13').Value = '  -0.90%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.95'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').Value = '63.770.95'
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('E16').Value = '  -5.39%  '
$ws.Range('D17').Value = '2.763.62'
$ws.Range('E17').Value = '  -0.89%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.14'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -2.45%  '
$ws.Range('E19').Value = '  -4.63%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '359.09'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -2.34%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.65'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -5.90%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('E23').Value = '  -7.95%  '
$ws.Range('E24').Value = '  -3.58%  '
$ws.Range('E25').Value = '  -3.35%  '
$ws.Range('E26').Value = '  -3.21%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').Value = '0.0₃0906'
$ws.Range('E28').Value = '  -6.89%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.37'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('E30').Value = '  -5.11%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.34'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +4.95%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '168.59'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -1.77%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.96'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  -5.02%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.22'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  -3.21%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.50'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -0.71%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('E38').Value = '  -2.15%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '348.25'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +1.76%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.33'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('E41').Value = '  -2.32%  '
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.10'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  -1.90%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.57'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -4.10%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0588'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -4.00%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '137.67'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  -0.88%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.628'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -3.73%  '
$ws.Range('E48').Value = '  -3.47%  '
$ws.Range('E51').Value = '  +0.10%  '
